# Applies the "Added editor and menu screen" logbook update:
#  1. The "Refactoring" bullet list gets a new explanatory sentence right
#     after "-Interface for encounters" (which itself gains a trailing
#     ": ").
#  2. The "Patterns implemented" bullet list gains a new "-Template Method"
#     entry.
#  3. The header byline is re-split (cosmetic run split only, the visible
#     text "Leander Suda, Florian Zolda" is unchanged).

$d = $word.ActiveDocument
$vt = [char]11   # manual line break (w:br) character used inside Word ranges

# --- 1. "Refactoring:" section -------------------------------------------
# "-Interface for encounters" -> "-Interface for encounters: " followed by
# a new line with the explanation of the encounter-interface unification.
$findText = "-Interface for encounters"
$replaceText = "-Interface for encounters: " + $vt + `
    "Encounters had a lot of code duplication and different ways to go. Unified the way encounters work with an interface."
$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $replaceText, 2) | Out-Null

# --- 2. "Patterns implemented" section ------------------------------------
# Add a new "-Template Method" bullet after the existing "-Dependency
# Injection" one.
$findText2 = "-Dependency Injection" + $vt
$replaceText2 = "-Dependency Injection" + $vt + "-Template Method"
$d.Content.Find.Execute($findText2, $true, $false, $false, $false, $false, `
    $true, 1, $false, $replaceText2, 2) | Out-Null

# --- 3. Header byline -------------------------------------------------------
# Split "Leander Suda, Florian Zolda" so "Zolda" lives in its own run
# (mirrors the spell-checker marking the surname as its own run); the
# visible text is unchanged.
$sec = $d.Sections.Item(1)
$hf = $sec.Headers.Item(1)
$headerText = $hf.Range.Text
$idx = $headerText.IndexOf("Zolda")
if ($idx -ge 0) {
    $wordStart = $hf.Range.Start + $idx
    $wordEnd = $wordStart + 5
    $zolda = $hf.Range.Duplicate()
    $zolda.Start = $wordStart
    $zolda.End = $wordEnd
    $zolda.Text = ""
    $insPoint = $hf.Range.Duplicate()
    $insPoint.Start = $wordStart
    $insPoint.End = $wordStart
    $insPoint.InsertAfter("Zolda")
}

Write-Output $d.Content.Text
